$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.042.80'
$ws.Range('E2').Value = '  -0.61%  '
$ws.Range('D3').Value = '1.640.07'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.30'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5050'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.48%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2574'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06433'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07721'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.72%  '
$ws.Range('D12').Value = '1.647.09'
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.242'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').Value = '1.866.74'
$ws.Range('E14').Value = '  -1.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5449'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').Value = '0.0₅7901'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.57'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '26.032.34'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '203.88'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.287'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.992'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.964'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.926'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +9.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.34'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.73%  '
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.71'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.735'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05049'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.97%  '
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.247'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.37%  '
$ws.Range('E33').Value = '  -0.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.541'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('E35').Value = '  -0.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8930'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.620'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5629'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('D39').Value = '1.145.46'
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.563'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.668'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('E44').Value = '  -3.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.69'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('D46').Value = '1.778.82'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '54.85'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05034'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.02%  '
